# Generate Report for Handback
# Update the handback/xliff timestamps recorded on the Overview, zh-cn and
# de-de sheets to reflect the newly generated report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on the Overview sheet
$wsOverview.Range("G2").Value = "2017-02-09 13:48:07"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (L2)
$wsZhCn.Range("H2").Value = "2017-02-09 13:47:47"
$wsZhCn.Range("L2").Value = "2017-02-09 13:48:32"

# de-de sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (L2)
$wsDeDe.Range("H2").Value = "2017-02-09 13:48:07"
$wsDeDe.Range("L2").Value = "2017-02-09 13:48:55"
